$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-78 down to 43-79.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new data record.
$ws.Range("A42").Value = 7
$ws.Range("B42").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C42").Value = "Ñuble"
$ws.Range("D42").Value = 44589
$ws.Range("E42").Value = 16
$ws.Range("F42").Value = 100112030
$ws.Range("G42").Value = "Poroto granado"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 60
$ws.Range("K42").Value = 23000
$ws.Range("L42").Value = 24000
$ws.Range("M42").Value = 23500
$ws.Range("N42").Value = "$/saco 25 kilos"
$ws.Range("O42").Value = "Provincia de Diguillín"
$ws.Range("P42").Value = 940
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"
